$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update names (column A)
$ws.Range("A2").Value = "Pedro"
$ws.Range("A3").Value = "Luiza"
$ws.Range("A4").Value = "Maria Rita"
$ws.Range("A5").Value = "Bruno"

# Update phone numbers (column B) - stored as text, not numbers
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "554812345678"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "5548123456789"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "554887654321"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "5548987654321"

# Update the active cell selection
$ws.Range("C3").Select()
